$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.281.51"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "3.537.73"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("D7").Value = "3.537.74"
$ws.Range("E7").Value = "  +0.64%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.478"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("E10").Value = "  -4.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.414"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "4.139.88"
$ws.Range("E14").Value = "  -3.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.58%  "
$ws.Range("D16").Value = "3.539.11"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").Value = "66.376.46"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "425.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("E24").Value = "  -1.40%  "
$ws.Range("D25").Value = "3.686.38"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.04%  "
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.24%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.159"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").Value = "3.526.20"
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("E39").Value = "  -4.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0862"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.892"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("E45").Value = "  -9.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("E47").Value = "  -7.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("E50").Value = "  -3.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.954"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.80%  "
